$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.695.10"
$ws.Range("E2").Value = "  -3.73%  "
$ws.Range("D3").Value = "3.045.79"
$ws.Range("E3").Value = "  -3.06%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'542.83"
$ws.Range("E5").Value = "  -4.35%  "
$ws.Range("D6").Value = "'133.67"
$ws.Range("E6").Value = "  -10.26%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.039.78"
$ws.Range("E8").Value = "  -2.95%  "
$ws.Range("D9").Value = "'0.487"
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("D10").Value = "'6.39"
$ws.Range("E10").Value = "  -10.51%  "
$ws.Range("D11").Value = "'0.154"
$ws.Range("E11").Value = "  -3.40%  "
$ws.Range("D12").Value = "'0.457"
$ws.Range("E12").Value = "  -2.14%  "
$ws.Range("D13").Value = "'34.63"
$ws.Range("E13").Value = "  -4.15%  "
$ws.Range("E14").Value = "  -4.69%  "
$ws.Range("D15").Value = "3.539.26"
$ws.Range("E15").Value = "  -3.04%  "
$ws.Range("D16").Value = "62.744.93"
$ws.Range("E16").Value = "  -3.66%  "
$ws.Range("E17").Value = "  -2.73%  "
$ws.Range("D18").Value = "3.045.29"
$ws.Range("E18").Value = "  -3.10%  "
$ws.Range("D19").Value = "'6.58"
$ws.Range("E19").Value = "  -2.96%  "
$ws.Range("D20").Value = "'478.44"
$ws.Range("E20").Value = "  -11.54%  "
$ws.Range("D21").Value = "'13.32"
$ws.Range("E21").Value = "  -4.45%  "
$ws.Range("D22").Value = "'0.693"
$ws.Range("E22").Value = "  -2.29%  "
$ws.Range("D23").Value = "'6.97"
$ws.Range("E23").Value = "  -6.91%  "
$ws.Range("D24").Value = "'77.06"
$ws.Range("E24").Value = "  -2.60%  "
$ws.Range("D25").Value = "'12.13"
$ws.Range("E25").Value = "  -5.94%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D27").Value = "'2.69"
$ws.Range("E27").Value = "  -4.66%  "
$ws.Range("D28").Value = "'8.20"
$ws.Range("E28").Value = "  -8.08%  "
$ws.Range("D30").Value = "'1.92"
$ws.Range("E30").Value = "  -10.65%  "
$ws.Range("D31").Value = "'26.07"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("E32").Value = "  -2.89%  "
$ws.Range("D33").Value = "'60.11"
$ws.Range("E33").Value = "  +13.38%  "
$ws.Range("D34").Value = "'2.48"
$ws.Range("E34").Value = "  -6.83%  "
$ws.Range("D35").Value = "'508.31"
$ws.Range("E35").Value = "  -8.49%  "
$ws.Range("D36").Value = "'5.90"
$ws.Range("E36").Value = "  -3.48%  "
$ws.Range("D37").Value = "'5.05"
$ws.Range("E37").Value = "  -7.19%  "
$ws.Range("D38").Value = "'0.0394"
$ws.Range("E38").Value = "  -12.10%  "
$ws.Range("D39").Value = "3.065.34"
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("D40").Value = "'0.0784"
$ws.Range("E40").Value = "  -4.76%  "
$ws.Range("E41").Value = "  -4.14%  "
$ws.Range("D42").Value = "'7.99"
$ws.Range("D43").Value = "'2.55"
$ws.Range("E43").Value = "  -12.44%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.250"
$ws.Range("E44").Value = "  -3.94%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "'2.01"
$ws.Range("E46").Value = "  -8.24%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'24.44"
$ws.Range("E47").Value = "  -2.85%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'119.07"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("E49").Value = "  -2.93%  "
$ws.Range("D50").Value = "0.0₃0491"
$ws.Range("E50").Value = "  -7.08%  "
$ws.Range("D51").Value = "'2.34"
$ws.Range("E51").Value = "  +59.45%  "
